$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.117.98'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '1.818.58'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '233.68'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '40.97'
$ws.Range('E8').Value = '  -3.72%  '
$ws.Range('E9').Value = '  +8.55%  '
$ws.Range('D10').Value = '0.0687'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('D11').Value = "'0.100"
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '2.081.07'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').Value = '1.822.17'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = '11.13'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '4.69'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.662'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').Value = '35.066.16'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '69.67'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('D19').Value = '0.0₃0793'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '239.98'
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').Value = '11.89'
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('D22').Value = '4.69'
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').Value = '172.91'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').Value = '7.89'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '17.54'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  +25.86%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = '4.08'
$ws.Range('E31').Value = '  +4.24%  '
$ws.Range('D32').Value = '3.338.02'
$ws.Range('E32').Value = '  +37.39%  '
$ws.Range('D33').Value = '0.0556'
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('D34').Value = "'4.00"
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = '1.79'
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('E36').Value = '  +7.86%  '
$ws.Range('D37').Value = '93.38'
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').Value = '0.685'
$ws.Range('E38').Value = '  +3.16%  '
$ws.Range('D39').Value = '0.0195'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.314.85'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '1.28'
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('D42').Value = '0.989'
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').Value = '14.74'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').Value = '2.33'
$ws.Range('E44').Value = '  -4.67%  '
$ws.Range('D45').Value = '2.46'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D47').Value = '6.37'
$ws.Range('E47').Value = '  +6.03%  '
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('D49').Value = '1.996.35'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '1.01'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '0.0649'
$ws.Range('E51').Value = '  +5.72%  '
